# Update column F (dSF) values on Sheet1 to reflect the repulled data / mean calculation.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$values = @{
    2  = 4
    3  = -4
    4  = 2
    5  = 1
    6  = 8
    8  = -3
    9  = 2
    10 = -1
    11 = -1
    12 = -1
    13 = -2
    14 = -2
    16 = -3
}

foreach ($row in $values.Keys) {
    $ws.Cells.Item($row, 6).Value = $values[$row]
}
